# Updating CPRA post-proc scripts to include MS River stations.
# Changing NOV-14 observation data from USACE to USGS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 observation: source USACE -> USGS, station id 01440 -> 291929089562600
$ws.Range("A12").Value = "USGS"
$ws.Range("B12").Value = "291929089562600"

# New MS River station column is wider than the others - best-fit column B
# (character width 9.14 renders to the stored OOXML width of 10).
$ws.Columns.Item(2).ColumnWidth = 9.14

# Move the active selection to reflect where the edit was made.
$ws.Range("F12").Select()
